# Auto-generated Excel COM-interop script
# Applies vm_pu.xlsx data update: column B (vm_pu at bus 0/slack) changes
# from 1.05 to 1.02 (380 kV case), and dependent bus voltage magnitudes
# (columns C-F, I-N) are recomputed accordingly for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers corresponding to spreadsheet columns B,C,D,E,F,I,J,K,L,M,N
$colNums = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# New values per row (row 2 .. row 25), one value per column in $colNums order
$data = @(
        @(1.02, 1.020187938760651, 1.024633669955879, 1.023805500138709, 1.018597342599396, 1.027414560603013, 1.025386462933698, 1.027461597882946, 1.026635858678833, 1.021443085277686, 1.026842629324535),
        @(1.02, 1.021323982055728, 1.025430354192305, 1.02488881807928, 1.020376162297516, 1.027575656475209, 1.026158611915067, 1.028065443163231, 1.0275253808059, 1.023025071810181, 1.027615874846024),
        @(1.02, 1.022058477354735, 1.025945150521496, 1.025589576627337, 1.021526318380475, 1.027678228622584, 1.02665715030907, 1.028454811698976, 1.028100156471565, 1.024047445146677, 1.02811512122173),
        @(1.02, 1.022367118116987, 1.026161401477948, 1.025884124811332, 1.02200964818018, 1.027720951002829, 1.026866475899014, 1.028618177825923, 1.028341601491726, 1.024476953683914, 1.028324744078021),
        @(1.02, 1.022418932066561, 1.026197701057487, 1.025933577834675, 1.022090790183684, 1.027728100892114, 1.026901607396962, 1.028645588726456, 1.028382130048764, 1.024549052879945, 1.028359925466726),
        @(1.02, 1.022062601978589, 1.025948040741922, 1.025593512591251, 1.021532777414131, 1.027678801048376, 1.026659948346986, 1.028456995878492, 1.028103383418532, 1.0240531854179, 1.02811792323318),
        @(1.02, 1.020571995940502, 1.02490306105021, 1.024171659455372, 1.019198686908034, 1.027469349181955, 1.025647641867842, 1.027665952442845, 1.026936644635275, 1.021977994831621, 1.027104179162723),
        @(1.02, 1.017940633590616, 1.023056179388646, 1.021664387967593, 1.015078697447956, 1.027087482582151, 1.023855368238419, 1.026261567901802, 1.024874450648534, 1.018311068447861, 1.025309360299169),
        @(1.02, 1.016183040366555, 1.02182116696458, 1.019991537887534, 1.012326694637258, 1.0268242883484, 1.022654710536264, 1.025318201698461, 1.02349532509878, 1.01585906365727, 1.024106997525361),
        @(1.02, 1.015421148261589, 1.021285486330238, 1.019266827364141, 1.011133641296108, 1.026708273687992, 1.022133407227809, 1.024908010140835, 1.022897092606284, 1.014795444569322, 1.023584953906413),
        @(1.02, 1.015138017479179, 1.021086372231881, 1.018997581428883, 1.010690263161104, 1.026664872225773, 1.021939557576338, 1.024755388579578, 1.022674720270383, 1.014400075357693, 1.023390828966202),
        @(1.02, 1.015198755938857, 1.0211290891706, 1.019055338183173, 1.010785379651117, 1.026674195952778, 1.021981148726106, 1.024788138130093, 1.02272242727091, 1.014484896857539, 1.023432479180174),
        @(1.02, 1.015397747257968, 1.021269030329615, 1.019244572573374, 1.011096996180705, 1.026704692404766, 1.02211738794301, 1.024895399670423, 1.022878714562568, 1.014762769280143, 1.023568911872392),
        @(1.02, 1.015520335027378, 1.021355234306054, 1.019361158598693, 1.011288963234535, 1.026723441387722, 1.022201300937905, 1.024961452790124, 1.022974986756385, 1.014933936412393, 1.023652944033365),
        @(1.02, 1.016233586581794, 1.021856698985821, 1.020039626813253, 1.012405842763617, 1.026831944630358, 1.022689277790365, 1.025345388646067, 1.023535005222472, 1.01592961194314, 1.024141613868927),
        @(1.02, 1.016680762351779, 1.022171009334034, 1.020465114890534, 1.013106043730899, 1.026899456709335, 1.022994993386448, 1.02558576313733, 1.023886003851146, 1.016553660849913, 1.024447763616223),
        @(1.02, 1.016941511351267, 1.022354253336047, 1.020713260786868, 1.013514323108861, 1.026938637618715, 1.023173176187349, 1.025725804845126, 1.024090633119825, 1.016917476876867, 1.024626199457139),
        @(1.02, 1.01703040637841, 1.022416719937253, 1.020797866424255, 1.013653513118616, 1.026951963769996, 1.023233908934942, 1.025773527565459, 1.024160389139441, 1.017041498333173, 1.0246870184522),
        @(1.02, 1.016632793035049, 1.022137295911968, 1.020419467595024, 1.013030932949815, 1.026892233759637, 1.022962207075274, 1.025559990273647, 1.023848355620494, 1.016486725123038, 1.024414930744726),
        @(1.02, 1.015339152894993, 1.021227824979949, 1.01918884933531, 1.011005239187331, 1.026695720478507, 1.022077274827472, 1.02486382095929, 1.022832696347669, 1.01468095097081, 1.023528741791629),
        @(1.02, 1.01452503532563, 1.020655202701902, 1.018414784360011, 1.00973029732949, 1.026570379924387, 1.021519641159386, 1.02442461729075, 1.022193171588926, 1.013543884411072, 1.022970316219773),
        @(1.02, 1.014956687174053, 1.020958837089969, 1.018825162728035, 1.010406296115425, 1.026636994672732, 1.021815371840187, 1.024657589649253, 1.02253228568674, 1.014146830036124, 1.023266466872061),
        @(1.02, 1.016654468554787, 1.022152529827865, 1.020440093748289, 1.013064872697811, 1.026895498110704, 1.022977022217935, 1.025571636430127, 1.023865367551352, 1.016516971055549, 1.024429766926589),
        @(1.02, 1.018621481643824, 1.023534300162317, 1.022312803142222, 1.016144709507044, 1.027187721466218, 1.024319728844246, 1.026625882655276, 1.025408330983802, 1.019260318638836, 1.025774380350321)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $colNums.Length; $j++) {
        $colNum = $colNums[$j]
        $ws.Cells.Item($rowNum, $colNum).Value = [double]$rowValues[$j]
    }
}
